$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 9; $r++) {
    $oldB = $ws.Cells.Item($r, 2).Value()

    # Copy the existing B cell (with its style) into C, D, E first so the
    # new cells inherit the same formatting/style as column B.
    $ws.Cells.Item($r, 2).Copy($ws.Cells.Item($r, 3))
    $ws.Cells.Item($r, 2).Copy($ws.Cells.Item($r, 4))
    $ws.Cells.Item($r, 2).Copy($ws.Cells.Item($r, 5))

    $ws.Cells.Item($r, 2).Value = 10044.0
    $ws.Cells.Item($r, 3).Value = $oldB
    $ws.Cells.Item($r, 4).Value = 20412.0
    $ws.Cells.Item($r, 5).Value = 32400.0
}
